$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the duplicate bold "Play 8 Dragons..." paragraph that sits
#    just before the final (italic) blurb paragraph -- it was moved up
#    to the top of the doc (see step 3) and relabeled as the meta
#    description, so this leftover copy near the bottom goes away.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$dupPara  = $lastPara.Previous(1)
$dupPara.Range.Delete()

# ------------------------------------------------------------------
# 2) Rewrite the final paragraph's text (the old meta-description blurb)
#    with the new image-prompt copy, preserving its italic run formatting.
#    Find.Execute collapses $blurbRng onto the matched text so setting
#    .Text only touches that run's text, leaving the <w:rPr><w:i/></w:rPr>
#    formatting (and the paragraph's leading empty run) untouched.
# ------------------------------------------------------------------
$oldBlurb = "Read our review of 8 Dragons, an online slot game with impeccable graphics and a chance to choose free spins and multipliers. Play for free today."
$newBlurb = "Create a feature image for ""8 Dragons"" that features a happy Maya warrior with glasses. The image should be in a cartoon style and should have a vibrant and eye-catching color scheme. The Maya warrior should be depicted holding a dragon in one hand and a pile of gold coins in the other, surrounded by Chinese-themed symbols such as lanterns and scrolls. In the background, you can add a colorful dragon or a temple to add to the overall theme of the game. The image should convey the excitement and adventure of playing ""8 Dragons"" and entice players to try their luck at this exciting slot game."
$blurbRng = $d.Content
$blurbRng.Find.Execute($oldBlurb, $true, $false, $false, $false, $false,
                        $true, 1, $false, "", 0)
$blurbRng.Text = $newBlurb

# ------------------------------------------------------------------
# 3) Insert a brand new "Meta description" paragraph right after the H1
#    title at the top of the document: an empty run, a bold
#    "Meta description" run, and a plain run with the colon-prefixed
#    blurb -- mirrors the paragraph that used to live at the bottom.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of 8 Dragons, an online slot game with impeccable graphics and a chance to choose free spins and multipliers. Play for free today.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaPara.Range.InsertXML($metaXml)

Write-Output "done"
